$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": new asesor "FARIAS CAICEDO GABRIELA PATRICIA"
# is inserted as a new row 3 (pushing the existing rows down by one), and
# column B is widened by one character.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Rows.Item(3).Insert()

$ws1.Range("A3").Value = "OFICINA-CATAECSA"
$ws1.Range("B3").Value = "FARIAS CAICEDO GABRIELA PATRICIA"
$cols1 = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")
foreach ($col in $cols1) {
    $ws1.Range($col + "3").Value = 0
}

# Update the "x de N" summary row (now row 6) to reflect the new headcount.
$ws1.Range("C6").Value = "0 de 4"
$ws1.Range("D6").Value = "1 de 4"
$ws1.Range("E6").Value = "0 de 4"
$ws1.Range("F6").Value = "0 de 4"
$ws1.Range("G6").Value = "0 de 4"
$ws1.Range("H6").Value = "0 de 4"
$ws1.Range("I6").Value = "0 de 4"
$ws1.Range("J6").Value = "0 de 4"
$ws1.Range("K6").Value = "0 de 4"
$ws1.Range("L6").Value = "0 de 4"
$ws1.Range("M6").Value = "0 de 4"
$ws1.Range("N6").Value = "0 de 4"
$ws1.Range("O6").Value = "0 de 4"
$ws1.Range("P6").Value = "0 de 4"
$ws1.Range("Q6").Value = "0 de 4"
$ws1.Range("R6").Value = "0 de 4"

$ws1.Columns.Item(2).ColumnWidth = 34 - 5/6

# ---------------------------------------------------------------------------
# Sheet "VENTA MENSUAL": same new asesor inserted as row 3, same column
# widening. Her monthly sales are all zero (she has not sold anything yet),
# the totals row (now row 6) keeps its original values unchanged.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(3).Insert()

$ws2.Range("A3").Value = "OFICINA-CATAECSA"
$ws2.Range("B3").Value = "FARIAS CAICEDO GABRIELA PATRICIA"
$ws2.Range("C3").Value = 0
$ws2.Range("D3").Value = 0
$ws2.Range("E3").Value = 0
$ws2.Range("F3").Value = 0
$ws2.Range("G3").Value = 0

$ws2.Columns.Item(2).ColumnWidth = 34 - 5/6
